{"js": "// Update the worksheet date and all 25 \"a\u00f7b=\" division-problem prompts.\n// Every \"old\" string below occurs exactly once in the document body, so a\n// plain text search + replace is unambiguous for each pair.\nconst replacements = [\n  [\"2024-03-23 Saturday\", \"2024-03-24 Sunday\"],\n  [\"46\u00f73=\", \"30\u00f78=\"],\n  [\"68\u00f77=\", \"88\u00f75=\"],\n  [\"86\u00f78=\", \"72\u00f73=\"],\n  [\"43\u00f78=\", \"32\u00f77=\"],\n  [\"19\u00f75=\", \"64\u00f75=\"],\n  [\"17\u00f76=\", \"12\u00f77=\"],\n  [\"21\u00f76=\", \"54\u00f74=\"],\n  [\"67\u00f79=\", \"32\u00f73=\"],\n  [\"88\u00f72=\", \"17\u00f79=\"],\n  [\"98\u00f74=\", \"19\u00f76=\"],\n  [\"26\u00f74=\", \"31\u00f78=\"],\n  [\"80\u00f78=\", \"14\u00f73=\"],\n  [\"10\u00f73=\", \"77\u00f75=\"],\n  [\"30\u00f74=\", \"91\u00f76=\"],\n  [\"23\u00f78=\", \"66\u00f74=\"],\n  [\"39\u00f77=\", \"47\u00f72=\"],\n  [\"13\u00f79=\", \"24\u00f76=\"],\n  [\"14\u00f75=\", \"78\u00f74=\"],\n  [\"81\u00f73=\", \"96\u00f78=\"],\n  [\"46\u00f74=\", \"72\u00f74=\"],\n  [\"14\u00f78=\", \"84\u00f79=\"],\n  [\"47\u00f79=\", \"12\u00f73=\"],\n  [\"61\u00f72=\", \"20\u00f79=\"],\n  [\"63\u00f75=\", \"94\u00f77=\"],\n  [\"21\u00f77=\", \"19\u00f74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items/text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-03-23 Saturday\", \"2024-03-24 Sunday\"),\n    @(\"46\u00f73=\", \"30\u00f78=\"),\n    @(\"68\u00f77=\", \"88\u00f75=\"),\n    @(\"86\u00f78=\", \"72\u00f73=\"),\n    @(\"43\u00f78=\", \"32\u00f77=\"),\n    @(\"19\u00f75=\", \"64\u00f75=\"),\n    @(\"17\u00f76=\", \"12\u00f77=\"),\n    @(\"21\u00f76=\", \"54\u00f74=\"),\n    @(\"67\u00f79=\", \"32\u00f73=\"),\n    @(\"88\u00f72=\", \"17\u00f79=\"),\n    @(\"98\u00f74=\", \"19\u00f76=\"),\n    @(\"26\u00f74=\", \"31\u00f78=\"),\n    @(\"80\u00f78=\", \"14\u00f73=\"),\n    @(\"10\u00f73=\", \"77\u00f75=\"),\n    @(\"30\u00f74=\", \"91\u00f76=\"),\n    @(\"23\u00f78=\", \"66\u00f74=\"),\n    @(\"39\u00f77=\", \"47\u00f72=\"),\n    @(\"13\u00f79=\", \"24\u00f76=\"),\n    @(\"14\u00f75=\", \"78\u00f74=\"),\n    @(\"81\u00f73=\", \"96\u00f78=\"),\n    @(\"46\u00f74=\", \"72\u00f74=\"),\n    @(\"14\u00f78=\", \"84\u00f79=\"),\n    @(\"47\u00f79=\", \"12\u00f73=\"),\n    @(\"61\u00f72=\", \"20\u00f79=\"),\n    @(\"63\u00f75=\", \"94\u00f77=\"),\n    @(\"21\u00f77=\", \"19\u00f74=\"),\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $found = $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n}\n\n"}
